# Apply "Made changes to Dashboard suite" edits to the Test Suite sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 = Dashboard_Page: Runmode (column C) changes from N to Y
$ws.Range("C6").Value = "Y"

# Row 8 = LeftHandPanel: Runmode (column C) changes from Y to N
$ws.Range("C8").Value = "N"

# Update the active cell selection to match the author's final cursor position
$ws.Range("J11").Select()
